# Update Physical Sales file to Exactly match Previous Royalty Run
# Previous royalty run had two separate physical sales files. This update
# merged those numbers so the Units (I column) reflect the combined totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RS Physical Sales Template")
$ws.Activate()

# --- Update the "Units" (column I) values to the merged totals ---------
$ws.Range("I2").Value = 142
$ws.Range("I7").Value = 148
$ws.Range("I9").Value = 37
$ws.Range("I12").Value = 245
$ws.Range("I13").Value = 198
$ws.Range("I16").Value = 51
$ws.Range("I17").Value = 91
$ws.Range("I24").Value = 11
$ws.Range("I25").Value = 13
$ws.Range("I26").Value = 13
$ws.Range("I27").Value = 14
$ws.Range("I28").Value = 18
$ws.Range("I29").Value = 51
$ws.Range("I40").Value = 103
$ws.Range("I42").Value = 48
$ws.Range("I44").Value = 40
$ws.Range("I45").Value = 85
$ws.Range("I50").Value = 103
$ws.Range("I52").Value = 12
$ws.Range("I53").Value = 14
$ws.Range("I54").Value = 14
$ws.Range("I55").Value = 48
$ws.Range("I56").Value = 16
$ws.Range("I57").Value = 11
$ws.Range("I58").Value = 11
$ws.Range("I60").Value = 11
$ws.Range("I61").Value = 12
$ws.Range("I62").Value = 45
$ws.Range("I63").Value = 455
$ws.Range("I64").Value = 161
$ws.Range("I66").Value = 220
$ws.Range("I67").Value = 112
$ws.Range("I69").Value = 15
$ws.Range("I70").Value = 16
$ws.Range("I71").Value = 51
$ws.Range("I72").Value = 10
$ws.Range("I73").Value = 153
$ws.Range("I74").Value = 176
$ws.Range("I75").Value = 224
$ws.Range("I76").Value = 41
$ws.Range("I77").Value = 348

# I2 picks up an integer number format as part of this edit.
$ws.Range("I2").NumberFormat = "0"

# --- Freeze panes at G2 (6 columns / 1 row frozen) and restore the
#     per-pane selections that Excel records once a sheet is split -------
$ws.Range("G2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true

# Pane 2 = top-right, Pane 3 = bottom-left, Pane 4 = bottom-right (active)
$win.Panes.Item(2).Activate()
$ws.Range("G1").Select()

$win.Panes.Item(3).Activate()
$ws.Range("A2").Select()

$win.Panes.Item(4).Activate()
$ws.Range("M7").Select()
